# Applies the "chore: update Sheets via scheduled runner" edit:
# refreshed market-price-derived values (currentAveragePrice*, Leve cost/profit
# columns H-N) across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets.
$wb = $excel.ActiveWorkbook

# ==== ALC ====
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(64, 8).Value2 = 7423.625  # H64: 7912.7144 -> 7423.625
$ws.Cells.Item(64, 9).Value2 = 6161.636  # I64: 6642 -> 6161.636
$ws.Cells.Item(64, 11).Value2 = 6161.636  # K64: 6642 -> 6161.636
$ws.Cells.Item(64, 13).Value2 = -5913.636  # M64: -6394 -> -5913.636
$ws.Cells.Item(67, 8).Value2 = 7423.625  # H67: 7912.7144 -> 7423.625
$ws.Cells.Item(67, 9).Value2 = 6161.636  # I67: 6642 -> 6161.636
$ws.Cells.Item(67, 11).Value2 = 6161.636  # K67: 6642 -> 6161.636
$ws.Cells.Item(67, 13).Value2 = -5303.636  # M67: -5784 -> -5303.636
$ws.Cells.Item(129, 8).Value2 = 253063.5  # H129: 337888.66 -> 253063.5
$ws.Cells.Item(129, 9).Value2 = 288755.84  # I129: 504286.5 -> 288755.84
$ws.Cells.Item(129, 10).Value2 = 3217  # J129: 5093 -> 3217
$ws.Cells.Item(129, 11).Value2 = 866267.52  # K129: 1512859.5 -> 866267.52
$ws.Cells.Item(129, 12).Value2 = 9651  # L129: 15279 -> 9651
$ws.Cells.Item(129, 13).Value2 = -861267.52  # M129: -1507859.5 -> -861267.52
$ws.Cells.Item(129, 14).Value2 = -19651  # N129: -25279 -> -19651
$ws.Cells.Item(137, 8).Value2 = 2326.415  # H137: 2371.577 -> 2326.415
$ws.Cells.Item(137, 9).Value2 = 1652.3438  # I137: 1706.3549 -> 1652.3438
$ws.Cells.Item(137, 11).Value2 = 4957.0314  # K137: 5119.0647 -> 4957.0314
$ws.Cells.Item(137, 13).Value2 = -2407.0314  # M137: -2569.0647 -> -2407.0314
$ws.Cells.Item(138, 8).Value2 = 4336.4  # H138: 4436.28 -> 4336.4
$ws.Cells.Item(138, 9).Value2 = 1500  # I138: 0 -> 1500
$ws.Cells.Item(138, 10).Value2 = 4454.5835  # J138: 4436.28 -> 4454.5835
$ws.Cells.Item(138, 11).Value2 = 4500  # K138: 0 -> 4500
$ws.Cells.Item(138, 12).Value2 = 13363.7505  # L138: 13308.84 -> 13363.7505
$ws.Cells.Item(138, 13).Value2 = 640  # M138: None -> 640
$ws.Cells.Item(138, 14).Value2 = -23643.7505  # N138: -23588.84 -> -23643.7505

# ==== ARM ====
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 12062.943  # H32: 12197.75 -> 12062.943
$ws.Cells.Item(32, 9).Value2 = 9314.421  # I32: 9435.946 -> 9314.421
$ws.Cells.Item(32, 11).Value2 = 9314.421  # K32: 9435.946 -> 9314.421
$ws.Cells.Item(32, 13).Value2 = -9027.421  # M32: -9148.946 -> -9027.421
$ws.Cells.Item(122, 8).Value2 = 4035.9583  # H122: 4284.727 -> 4035.9583
$ws.Cells.Item(122, 9).Value2 = 2997.8125  # I122: 3240.4285 -> 2997.8125
$ws.Cells.Item(122, 11).Value2 = 8993.4375  # K122: 9721.2855 -> 8993.4375
$ws.Cells.Item(122, 13).Value2 = -6543.4375  # M122: -7271.2855 -> -6543.4375

# ==== BSM ====
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value2 = 1807.4584  # H20: 1856.0435 -> 1807.4584
$ws.Cells.Item(20, 9).Value2 = 1457.6  # I20: 1512.4286 -> 1457.6
$ws.Cells.Item(20, 11).Value2 = 1457.6  # K20: 1512.4286 -> 1457.6
$ws.Cells.Item(20, 13).Value2 = -1210.6  # M20: -1265.4286 -> -1210.6
$ws.Cells.Item(102, 8).Value2 = 26328.846  # H102: 26872.584 -> 26328.846
$ws.Cells.Item(102, 9).Value2 = 25189.584  # I102: 25679.182 -> 25189.584
$ws.Cells.Item(102, 11).Value2 = 25189.584  # K102: 25679.182 -> 25189.584
$ws.Cells.Item(102, 13).Value2 = -21944.584  # M102: -22434.182 -> -21944.584
$ws.Cells.Item(105, 8).Value2 = 2901.2559  # H105: 3081.65 -> 2901.2559
$ws.Cells.Item(105, 9).Value2 = 2115.7097  # I105: 2289.25 -> 2115.7097
$ws.Cells.Item(105, 11).Value2 = 2115.7097  # K105: 2289.25 -> 2115.7097
$ws.Cells.Item(105, 13).Value2 = -368.7096999999999  # M105: -542.25 -> -368.7096999999999
$ws.Cells.Item(107, 8).Value2 = 807.26666  # H107: 824.2143 -> 807.26666
$ws.Cells.Item(107, 9).Value2 = 800.8333  # I107: 821.8182 -> 800.8333
$ws.Cells.Item(107, 11).Value2 = 800.8333  # K107: 821.8182 -> 800.8333
$ws.Cells.Item(107, 13).Value2 = 1119.1667  # M107: 1098.1818 -> 1119.1667
$ws.Cells.Item(137, 8).Value2 = 69989.414  # H137: 69989.47 -> 69989.414
$ws.Cells.Item(137, 10).Value2 = 69989.414  # J137: 69989.47 -> 69989.414
$ws.Cells.Item(137, 12).Value2 = 69989.414  # L137: 69989.47 -> 69989.414
$ws.Cells.Item(137, 14).Value2 = -80189.414  # N137: -80189.47 -> -80189.414
$ws.Cells.Item(141, 8).Value2 = 59985  # H141: 59990.332 -> 59985
$ws.Cells.Item(141, 10).Value2 = 59970  # J141: 59971 -> 59970
$ws.Cells.Item(141, 12).Value2 = 59970  # L141: 59971 -> 59970
$ws.Cells.Item(141, 14).Value2 = -70330  # N141: -70331 -> -70330

# ==== CRP ====
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value2 = 576.2857  # H16: 666.75 -> 576.2857
$ws.Cells.Item(16, 9).Value2 = 488.8  # I16: 624 -> 488.8
$ws.Cells.Item(16, 11).Value2 = 488.8  # K16: 624 -> 488.8
$ws.Cells.Item(16, 13).Value2 = -201.8  # M16: -337 -> -201.8
$ws.Cells.Item(62, 8).Value2 = 5204.875  # H62: 3716.6924 -> 5204.875
$ws.Cells.Item(62, 9).Value2 = 6027.3335  # I62: 4514.222 -> 6027.3335
$ws.Cells.Item(62, 10).Value2 = 2737.5  # J62: 1922.25 -> 2737.5
$ws.Cells.Item(62, 11).Value2 = 6027.3335  # K62: 4514.222 -> 6027.3335
$ws.Cells.Item(62, 12).Value2 = 2737.5  # L62: 1922.25 -> 2737.5
$ws.Cells.Item(62, 13).Value2 = -5403.3335  # M62: -3890.222 -> -5403.3335
$ws.Cells.Item(62, 14).Value2 = -3985.5  # N62: -3170.25 -> -3985.5
$ws.Cells.Item(65, 8).Value2 = 5204.875  # H65: 3716.6924 -> 5204.875
$ws.Cells.Item(65, 9).Value2 = 6027.3335  # I65: 4514.222 -> 6027.3335
$ws.Cells.Item(65, 10).Value2 = 2737.5  # J65: 1922.25 -> 2737.5
$ws.Cells.Item(65, 11).Value2 = 30136.6675  # K65: 22571.11 -> 30136.6675
$ws.Cells.Item(65, 12).Value2 = 13687.5  # L65: 9611.25 -> 13687.5
$ws.Cells.Item(65, 13).Value2 = -27016.6675  # M65: -19451.11 -> -27016.6675
$ws.Cells.Item(65, 14).Value2 = -19927.5  # N65: -15851.25 -> -19927.5
$ws.Cells.Item(113, 8).Value2 = 576.2857  # H113: 666.75 -> 576.2857
$ws.Cells.Item(113, 9).Value2 = 488.8  # I113: 624 -> 488.8
$ws.Cells.Item(113, 11).Value2 = 488.8  # K113: 624 -> 488.8
$ws.Cells.Item(113, 13).Value2 = 1681.2  # M113: 1546 -> 1681.2
$ws.Cells.Item(132, 8).Value2 = 3097.0908  # H132: 3316.95 -> 3097.0908
$ws.Cells.Item(132, 9).Value2 = 3082.5334  # I132: 3456.923 -> 3082.5334
$ws.Cells.Item(132, 10).Value2 = 3128.2856  # J132: 3057 -> 3128.2856
$ws.Cells.Item(132, 11).Value2 = 9247.600199999999  # K132: 10370.769 -> 9247.600199999999
$ws.Cells.Item(132, 12).Value2 = 9384.856800000001  # L132: 9171 -> 9384.856800000001
$ws.Cells.Item(132, 13).Value2 = -6717.600199999999  # M132: -7840.769 -> -6717.600199999999
$ws.Cells.Item(132, 14).Value2 = -14444.8568  # N132: -14231 -> -14444.8568

# ==== CUL ====
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value2 = 368.7143  # H2: 337 -> 368.7143
$ws.Cells.Item(2, 10).Value2 = 458.2  # J2: 391.85715 -> 458.2
$ws.Cells.Item(2, 12).Value2 = 2749.2  # L2: 2351.1429 -> 2749.2
$ws.Cells.Item(2, 14).Value2 = -2975.2  # N2: -2577.1429 -> -2975.2
$ws.Cells.Item(7, 8).Value2 = 620.0952  # H7: 564.9583 -> 620.0952
$ws.Cells.Item(7, 9).Value2 = 316.75  # I7: 298.22223 -> 316.75
$ws.Cells.Item(7, 10).Value2 = 1590.8  # J7: 1365.1666 -> 1590.8
$ws.Cells.Item(7, 11).Value2 = 950.25  # K7: 894.66669 -> 950.25
$ws.Cells.Item(7, 12).Value2 = 4772.4  # L7: 4095.4998 -> 4772.4
$ws.Cells.Item(7, 13).Value2 = -838.25  # M7: -782.66669 -> -838.25
$ws.Cells.Item(7, 14).Value2 = -4996.4  # N7: -4319.4998 -> -4996.4
$ws.Cells.Item(34, 8).Value2 = 5014.222  # H34: 4205 -> 5014.222
$ws.Cells.Item(34, 9).Value2 = 892.6667  # I34: 761 -> 892.6667
$ws.Cells.Item(34, 11).Value2 = 2678.0001  # K34: 2283 -> 2678.0001
$ws.Cells.Item(34, 13).Value2 = -2594.0001  # M34: -2199 -> -2594.0001
$ws.Cells.Item(131, 8).Value2 = 3617.4  # H131: 3665.2917 -> 3617.4
$ws.Cells.Item(131, 10).Value2 = 4173  # J131: 4328 -> 4173
$ws.Cells.Item(131, 12).Value2 = 12519  # L131: 12984 -> 12519
$ws.Cells.Item(131, 14).Value2 = -22599  # N131: -23064 -> -22599
$ws.Cells.Item(138, 8).Value2 = 6240.1875  # H138: 6895.3335 -> 6240.1875
$ws.Cells.Item(138, 9).Value2 = 4896.8887  # I138: 5394.6 -> 4896.8887
$ws.Cells.Item(138, 11).Value2 = 14690.6661  # K138: 16183.8 -> 14690.6661
$ws.Cells.Item(138, 13).Value2 = -9550.666100000002  # M138: -11043.8 -> -9550.666100000002

# ==== GSM ====
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value2 = 10870205  # H97: 11905455 -> 10870205
$ws.Cells.Item(97, 9).Value2 = 658.86664  # I97: 701.7143 -> 658.86664
$ws.Cells.Item(97, 10).Value2 = 31250602  # J97: 35714960 -> 31250602
$ws.Cells.Item(97, 11).Value2 = 658.86664  # K97: 701.7143 -> 658.86664
$ws.Cells.Item(97, 12).Value2 = 31250602  # L97: 35714960 -> 31250602
$ws.Cells.Item(97, 13).Value2 = -162.86664  # M97: -205.7143 -> -162.86664
$ws.Cells.Item(97, 14).Value2 = -31251594  # N97: -35715952 -> -31251594
$ws.Cells.Item(113, 8).Value2 = 5486.7  # H113: 6991.4287 -> 5486.7
$ws.Cells.Item(113, 9).Value2 = 4415.4287  # I113: 5390 -> 4415.4287
$ws.Cells.Item(113, 10).Value2 = 7986.3335  # J113: 10995 -> 7986.3335
$ws.Cells.Item(113, 11).Value2 = 4415.4287  # K113: 5390 -> 4415.4287
$ws.Cells.Item(113, 12).Value2 = 7986.3335  # L113: 10995 -> 7986.3335
$ws.Cells.Item(113, 13).Value2 = -2245.4287  # M113: -3220 -> -2245.4287
$ws.Cells.Item(113, 14).Value2 = -12326.3335  # N113: -15335 -> -12326.3335
$ws.Cells.Item(132, 8).Value2 = 3297.2307  # H132: 3137.1428 -> 3297.2307
$ws.Cells.Item(132, 9).Value2 = 3357.6843  # I132: 3137.1428 -> 3357.6843
$ws.Cells.Item(132, 10).Value2 = 1000  # J132: 0 -> 1000
$ws.Cells.Item(132, 11).Value2 = 10073.0529  # K132: 9411.428400000001 -> 10073.0529
$ws.Cells.Item(132, 12).Value2 = 3000  # L132: 0 -> 3000
$ws.Cells.Item(132, 13).Value2 = -7543.052899999999  # M132: -6881.428400000001 -> -7543.052899999999
$ws.Cells.Item(132, 14).Value2 = -8060  # N132: None -> -8060
$ws.Cells.Item(133, 8).Value2 = 71944.16  # H133: 69662.42999999999 -> 71944.16
$ws.Cells.Item(133, 10).Value2 = 71944.16  # J133: 69662.42999999999 -> 71944.16
$ws.Cells.Item(133, 12).Value2 = 71944.16  # L133: 69662.42999999999 -> 71944.16
$ws.Cells.Item(133, 14).Value2 = -82064.16  # N133: -79782.42999999999 -> -82064.16

# ==== LTW ====
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value2 = 8114.4473  # H46: 7924.3335 -> 8114.4473
$ws.Cells.Item(46, 10).Value2 = 9703.226000000001  # J46: 9421.875 -> 9703.226000000001
$ws.Cells.Item(46, 12).Value2 = 9703.226000000001  # L46: 9421.875 -> 9703.226000000001
$ws.Cells.Item(46, 14).Value2 = -10079.226  # N46: -9797.875 -> -10079.226
$ws.Cells.Item(68, 8).Value2 = 5802.15  # H68: 5370.4585 -> 5802.15
$ws.Cells.Item(68, 9).Value2 = 5172.769  # I68: 4730.9375 -> 5172.769
$ws.Cells.Item(68, 10).Value2 = 6971  # J68: 6649.5 -> 6971
$ws.Cells.Item(68, 11).Value2 = 5172.769  # K68: 4730.9375 -> 5172.769
$ws.Cells.Item(68, 12).Value2 = 6971  # L68: 6649.5 -> 6971
$ws.Cells.Item(68, 13).Value2 = -4423.769  # M68: -3981.9375 -> -4423.769
$ws.Cells.Item(68, 14).Value2 = -8469  # N68: -8147.5 -> -8469
$ws.Cells.Item(71, 8).Value2 = 5802.15  # H71: 5370.4585 -> 5802.15
$ws.Cells.Item(71, 9).Value2 = 5172.769  # I71: 4730.9375 -> 5172.769
$ws.Cells.Item(71, 10).Value2 = 6971  # J71: 6649.5 -> 6971
$ws.Cells.Item(71, 11).Value2 = 25863.845  # K71: 23654.6875 -> 25863.845
$ws.Cells.Item(71, 12).Value2 = 34855  # L71: 33247.5 -> 34855
$ws.Cells.Item(71, 13).Value2 = -22119.845  # M71: -19910.6875 -> -22119.845
$ws.Cells.Item(71, 14).Value2 = -42343  # N71: -40735.5 -> -42343
$ws.Cells.Item(122, 8).Value2 = 5920.467  # H122: 7098.467 -> 5920.467
$ws.Cells.Item(122, 9).Value2 = 3655.7778  # I122: 5619.1113 -> 3655.7778
$ws.Cells.Item(122, 11).Value2 = 10967.3334  # K122: 16857.3339 -> 10967.3334
$ws.Cells.Item(122, 13).Value2 = -8517.3334  # M122: -14407.3339 -> -8517.3334
$ws.Cells.Item(127, 8).Value2 = 66905  # H127: 70000 -> 66905
$ws.Cells.Item(127, 10).Value2 = 66905  # J127: 70000 -> 66905
$ws.Cells.Item(127, 12).Value2 = 66905  # L127: 70000 -> 66905
$ws.Cells.Item(127, 14).Value2 = -76825  # N127: -79920 -> -76825
$ws.Cells.Item(132, 8).Value2 = 4640  # H132: 4715.394 -> 4640
$ws.Cells.Item(132, 9).Value2 = 3963.6538  # I132: 4036.12 -> 3963.6538
$ws.Cells.Item(132, 11).Value2 = 11890.9614  # K132: 12108.36 -> 11890.9614
$ws.Cells.Item(132, 13).Value2 = -9360.9614  # M132: -9578.360000000001 -> -9360.9614

# ==== WVR ====
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value2 = 11114056  # H62: 10103901 -> 11114056
$ws.Cells.Item(62, 9).Value2 = 2483.3333  # I62: 2450 -> 2483.3333
$ws.Cells.Item(62, 11).Value2 = 2483.3333  # K62: 2450 -> 2483.3333
$ws.Cells.Item(62, 13).Value2 = -1859.3333  # M62: -1826 -> -1859.3333
$ws.Cells.Item(65, 8).Value2 = 11114056  # H65: 10103901 -> 11114056
$ws.Cells.Item(65, 9).Value2 = 2483.3333  # I65: 2450 -> 2483.3333
$ws.Cells.Item(65, 11).Value2 = 12416.6665  # K65: 12250 -> 12416.6665
$ws.Cells.Item(65, 13).Value2 = -9296.666499999999  # M65: -9130 -> -9296.666499999999
$ws.Cells.Item(68, 8).Value2 = 0  # H68: 50000 -> 0
$ws.Cells.Item(68, 10).Value2 = 0  # J68: 50000 -> 0
$ws.Cells.Item(68, 12).Value2 = 0  # L68: 50000 -> 0
$ws.Cells.Item(68, 14).ClearContents()  # N68: -51622 -> (removed)
$ws.Cells.Item(71, 8).Value2 = 0  # H71: 50000 -> 0
$ws.Cells.Item(71, 10).Value2 = 0  # J71: 50000 -> 0
$ws.Cells.Item(71, 12).Value2 = 0  # L71: 150000 -> 0
$ws.Cells.Item(71, 14).ClearContents()  # N71: -158112 -> (removed)
$ws.Cells.Item(107, 8).Value2 = 494.42856  # H107: 633.44446 -> 494.42856
$ws.Cells.Item(107, 9).Value2 = 426.83334  # I107: 557.2857 -> 426.83334
$ws.Cells.Item(107, 11).Value2 = 1280.50002  # K107: 1671.8571 -> 1280.50002
$ws.Cells.Item(107, 13).Value2 = 639.4999800000001  # M107: 248.1428999999998 -> 639.4999800000001
$ws.Cells.Item(113, 8).Value2 = 568.11536  # H113: 591.7083 -> 568.11536
$ws.Cells.Item(113, 9).Value2 = 604.3333  # I113: 644.25 -> 604.3333
$ws.Cells.Item(113, 11).Value2 = 1812.9999  # K113: 1932.75 -> 1812.9999
$ws.Cells.Item(113, 13).Value2 = 357.0001  # M113: 237.25 -> 357.0001
$ws.Cells.Item(126, 8).Value2 = 2077.25  # H126: 2166.8572 -> 2077.25
$ws.Cells.Item(126, 9).Value2 = 1659.7142  # I126: 1694.6666 -> 1659.7142
$ws.Cells.Item(126, 11).Value2 = 4979.142599999999  # K126: 5083.9998 -> 4979.142599999999
$ws.Cells.Item(126, 13).Value2 = -2509.142599999999  # M126: -2613.9998 -> -2509.142599999999
